$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The city list (column A) was reshuffled in the shared-string table while
# the frequency counts (column B) stayed attached to their original rows.
# Net effect: only the text of certain A-cells changes. Apply those
# directly.
$ws.Range("A9").Value = 'Graz'
$ws.Range("A10").Value = 'Leipzig'
$ws.Range("A25").Value = 'Marburg'
$ws.Range("A26").Value = 'Tübingen'
$ws.Range("A27").Value = 'Mannheim'
$ws.Range("A28").Value = 'Passau'
$ws.Range("A29").Value = 'Bochum'
$ws.Range("A30").Value = 'Jena'
$ws.Range("A34").Value = 'Moskau'
$ws.Range("A35").Value = 'Bamberg'
$ws.Range("A36").Value = 'Wuppertal'
$ws.Range("A37").Value = 'Krems'
$ws.Range("A38").Value = 'Münster'
$ws.Range("A39").Value = 'Paris'
$ws.Range("A81").Value = 'Amsterdam'
$ws.Range("A82").Value = 'Antwerpen'
$ws.Range("A83").Value = 'Mailand'
$ws.Range("A84").Value = 'Kopenhagen'
$ws.Range("A87").Value = 'Athens'
$ws.Range("A93").Value = 'London'
$ws.Range("A97").Value = 'Glasgow'
$ws.Range("A104").Value = 'Cork'
